$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates ---
$ws.Range("A8").Value = "Volume 30   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/20/2023  Through  11/26/2023"

# --- Weekly crime statistics table updates (rows 14-30) ---
$ws.Range("D14").Value = 1
$ws.Range("G14").Value = 5
$ws.Range("J14").Value = 10
$ws.Range("K14").Value = -80
$ws.Range("N14").Value = -90.909090909090
$ws.Range("D15").Value = '0'
$ws.Range("E15").Value = '***.*'
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("M15").Value = 22.222222222222
$ws.Range("A16").Value = 'Robbery'
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("I16").Value = 156
$ws.Range("J16").Value = 165
$ws.Range("K16").Value = -5.454545454545
$ws.Range("L16").Value = 20
$ws.Range("M16").Value = -49.677419354838
$ws.Range("N16").Value = -83.439490445859
$ws.Range("A17").Value = 'Fel. Assault'
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = -57.894736842105
$ws.Range("I17").Value = 393
$ws.Range("J17").Value = 390
$ws.Range("K17").Value = 0.769230769230
$ws.Range("L17").Value = 5.645161290322
$ws.Range("M17").Value = 51.737451737451
$ws.Range("N17").Value = 3.149606299212
$ws.Range("A18").Value = 'Burglary'
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 21
$ws.Range("H18").Value = 10.526315789473
$ws.Range("I18").Value = 213
$ws.Range("J18").Value = 193
$ws.Range("K18").Value = 10.362694300518
$ws.Range("L18").Value = 51.063829787234
$ws.Range("M18").Value = -39.316239316239
$ws.Range("N18").Value = -85.733422638981
$ws.Range("A19").Value = 'Gr. Larceny'
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -9.090909090909
$ws.Range("F19").Value = 46
$ws.Range("H19").Value = -9.803921568627
$ws.Range("I19").Value = 626
$ws.Range("J19").Value = 572
$ws.Range("K19").Value = 9.440559440559
$ws.Range("L19").Value = 44.907407407407
$ws.Range("M19").Value = 53.431372549019
$ws.Range("N19").Value = 15.711645101663
$ws.Range("A20").Value = 'G.L.A.'
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = -17.391304347826
$ws.Range("I20").Value = 312
$ws.Range("J20").Value = 244
$ws.Range("K20").Value = 27.868852459016
$ws.Range("L20").Value = 124.460431654676
$ws.Range("M20").Value = -11.864406779661
$ws.Range("N20").Value = -90.076335877862
$ws.Range("A21").Value = 'TOTAL'
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -24.137931034482
$ws.Range("F21").Value = 118
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = -19.178082191780
$ws.Range("I21").Value = 1724
$ws.Range("J21").Value = 1587
$ws.Range("K21").Value = 8.632640201638
$ws.Range("L21").Value = 37.261146496815
$ws.Range("M21").Value = 0.583430571761
$ws.Range("N21").Value = -73.715505412410
$ws.Range("A22").Value = 'Transit'
$ws.Range("E22").Value = '***.*'
$ws.Range("H22").Value = '***.*'
$ws.Range("K22").Value = '***.*'
$ws.Range("L22").Value = '***.*'
$ws.Range("M22").Value = '***.*'
$ws.Range("N22").Value = '***.*'
$ws.Range("E23").Value = '***.*'
$ws.Range("H23").Value = '***.*'
$ws.Range("K23").Value = '***.*'
$ws.Range("L23").Value = '***.*'
$ws.Range("M23").Value = '***.*'
$ws.Range("N23").Value = '***.*'
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = 23.076923076923
$ws.Range("F24").Value = 79
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = -24.761904761904
$ws.Range("I24").Value = 1129
$ws.Range("J24").Value = 1375
$ws.Range("K24").Value = -17.890909090909
$ws.Range("L24").Value = 17.726798748696
$ws.Range("M24").Value = 49.734748010610
$ws.Range("N24").Value = '***.*'
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -7.142857142857
$ws.Range("F25").Value = 63
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 65.789473684210
$ws.Range("I25").Value = 638
$ws.Range("J25").Value = 596
$ws.Range("K25").Value = 7.046979865771
$ws.Range("L25").Value = 37.796976241900
$ws.Range("M25").Value = 16.849816849816
$ws.Range("N25").Value = '***.*'
$ws.Range("D26").Value = '0'
$ws.Range("E26").Value = '***.*'
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 200
$ws.Range("M26").Value = '***.*'
$ws.Range("N26").Value = '***.*'
$ws.Range("C27").Value = '0'
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 57
$ws.Range("K27").Value = -22.807017543859
$ws.Range("M27").Value = '***.*'
$ws.Range("N27").Value = '***.*'
$ws.Range("D28").Value = 1
$ws.Range("G28").Value = 3
$ws.Range("J28").Value = 25
$ws.Range("K28").Value = -68
$ws.Range("M28").Value = -82.222222222222
$ws.Range("N28").Value = -89.189189189189
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 20
$ws.Range("K29").Value = -65
$ws.Range("M29").Value = -76.666666666666
$ws.Range("N29").Value = -89.393939393939
$ws.Range("E30").Value = '***.*'
$ws.Range("H30").Value = '***.*'
$ws.Range("M30").Value = '***.*'
$ws.Range("N30").Value = '***.*'
$ws.Range("A38").Value = 'Robbery'
$ws.Range("A39").Value = 'Fel. Assault'
$ws.Range("A40").Value = 'Burglary'
$ws.Range("A41").Value = 'Gr. Larceny'
$ws.Range("A42").Value = 'G.L.A.'
$ws.Range("A43").Value = 'TOTAL'
